# Weekly update: insert a new price record for Maracuyá at Vega Modelo de
# Temuco (row 40), shifting all the existing records below it down by one
# row (old row 40 becomes row 41, ..., old row 116 becomes row 117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40; this shifts rows 40..116 down to 41..117
# and keeps their values/formatting intact.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with this week's record.
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = 'Vega Modelo de Temuco'
$ws.Range("C40").Value = 'La Araucanía'
$ws.Range("D40").Value = '2023-08-24'
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = 'Fruta'
$ws.Range("G40").Value = 100108
$ws.Range("H40").Value = 'Tropicales y subtropicales'
$ws.Range("I40").Value = 100108003
$ws.Range("J40").Value = 'Maracuyá'
$ws.Range("K40").Value = 'Sin especificar'
$ws.Range("L40").Value = 'Primera'
$ws.Range("M40").Value = 100
$ws.Range("N40").Value = 38000
$ws.Range("O40").Value = 38000
$ws.Range("P40").Value = 38000
$ws.Range("Q40").Value = '$/caja 18 kilos'
$ws.Range("R40").Value = 'Región de Arica y Parinacota'
$ws.Range("S40").Value = 2111
$ws.Range("T40").Value = 18
